$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price column (D) - force text format to preserve exact formatting (trailing zeros, dot-grouped numbers)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.399.43'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.452.21'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.30'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '179.37'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.588'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '48.61'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '684.65'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.002.05'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.70'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '69.499.68'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.449.67'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.91'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.39'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.911'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '17.07'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '101.11'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.71'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.74'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '33.78'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.82'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.91'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.75'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '563.09'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.08'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '58.20'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.633.96'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '35.12'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0₃0744'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0426'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.337'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.00'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '131.50'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.70'

# Update Volume(1h) column (E)
$ws.Range("E2").Value = '  +0.43%  '
$ws.Range("E3").Value = '  +2.48%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("E5").Value = '  -0.43%  '
$ws.Range("E6").Value = '  +0.07%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  +7.19%  '
$ws.Range("E10").Value = '  +0.45%  '
$ws.Range("E11").Value = '  +0.79%  '
$ws.Range("E12").Value = '  +2.38%  '
$ws.Range("E13").Value = '  -1.89%  '
$ws.Range("E14").Value = '  +1.99%  '
$ws.Range("E15").Value = '  +2.14%  '
$ws.Range("E16").Value = '  +0.57%  '
$ws.Range("E17").Value = '  +2.09%  '
$ws.Range("E19").Value = '  +1.68%  '
$ws.Range("E20").Value = '  +1.04%  '
$ws.Range("E21").Value = '  +1.33%  '
$ws.Range("E22").Value = '  -2.09%  '
$ws.Range("E23").Value = '  -0.24%  '
$ws.Range("E24").Value = '  -0.38%  '
$ws.Range("E25").Value = '  -0.26%  '
$ws.Range("E26").Value = '  -0.13%  '
$ws.Range("E27").Value = '  +1.55%  '
$ws.Range("E28").Value = '  +0.91%  '
$ws.Range("E29").Value = '  +2.34%  '
$ws.Range("E30").Value = '  -1.90%  '
$ws.Range("E31").Value = '  +5.81%  '
$ws.Range("E32").Value = '  +1.46%  '
$ws.Range("E33").Value = '  -0.53%  '
$ws.Range("E34").Value = '  -0.61%  '
$ws.Range("E35").Value = '  +0.16%  '
$ws.Range("E36").Value = '  +0.11%  '
$ws.Range("E37").Value = '  -2.07%  '
$ws.Range("E38").Value = '  -1.50%  '
$ws.Range("E39").Value = '  +0.83%  '
$ws.Range("E40").Value = '  +9.24%  '
$ws.Range("E41").Value = '  +2.47%  '
$ws.Range("E43").Value = '  +2.99%  '
$ws.Range("E44").Value = '  +2.14%  '
$ws.Range("E45").Value = '  -1.02%  '
$ws.Range("E46").Value = '  +0.93%  '
$ws.Range("E47").Value = '  -0.05%  '
$ws.Range("E48").Value = '  +5.17%  '
$ws.Range("E49").Value = '  -0.14%  '
$ws.Range("E50").Value = '  -0.26%  '
$ws.Range("E51").Value = '  +2.32%  '
